# feat: Add SFX of Game Start & End resource file
#
# Adds four new SFX rows (GameOver / LevelClear / LevelStartBeep / LevelStartCount)
# for the "GameplayHUD" / "GameScreen" team+situation to the SFX sheet, keeping
# the sheet sorted alphabetically by the generated FileName (column G), which
# pushes the pre-existing "MenuSelect" row down to row 6.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("SFX")

# Move the existing MainMenu/MenuSelect row down to row 6 to make room above it
# (rows are ordered alphabetically by the generated G column value, and
# "GameplayHUD_..." sorts before "MainMenu_...").
$ws.Range("A6").Value = "MainMenu"
$ws.Range("B6").Value = "TitleScreen"
$ws.Range("C6").Value = "MenuSelect"
$ws.Range("D6").Value = "BlipLow.wav"
$ws.Range("E6").Value = "S"
$ws.Range("F6").Value = "O"
$ws.Range("G6").Formula = "=_xlfn.CONCAT(""SFX_"",A6,""_"",B6,""_"",C6,""_"",D6)"

# Row 2: GameOver
$ws.Range("A2").Value = "GameplayHUD"
$ws.Range("B2").Value = "GameScreen"
$ws.Range("C2").Value = "GameOver"
$ws.Range("D2").Value = "DescendingScales2.wav"
$ws.Range("E2").Value = "S"
$ws.Range("F2").Value = "O"
$ws.Range("G2").Formula = "=_xlfn.CONCAT(""SFX_"",A2,""_"",B2,""_"",C2,""_"",D2)"

# Row 3: LevelClear
$ws.Range("A3").Value = "GameplayHUD"
$ws.Range("B3").Value = "GameScreen"
$ws.Range("C3").Value = "LevelClear"
$ws.Range("D3").Value = "AscendingScales2.wav"
$ws.Range("E3").Value = "S"
$ws.Range("F3").Value = "O"
$ws.Range("G3").Formula = "=_xlfn.CONCAT(""SFX_"",A3,""_"",B3,""_"",C3,""_"",D3)"

# Row 4: LevelStartBeep
$ws.Range("A4").Value = "GameplayHUD"
$ws.Range("B4").Value = "GameScreen"
$ws.Range("C4").Value = "LevelStartBeep"
$ws.Range("D4").Value = "BlipHiLong.wav"
$ws.Range("E4").Value = "S"
$ws.Range("F4").Value = "O"
$ws.Range("G4").Formula = "=_xlfn.CONCAT(""SFX_"",A4,""_"",B4,""_"",C4,""_"",D4)"

# Row 5: LevelStartCount
$ws.Range("A5").Value = "GameplayHUD"
$ws.Range("B5").Value = "GameScreen"
$ws.Range("C5").Value = "LevelStartCount"
$ws.Range("D5").Value = "BlipHiShort.wav"
$ws.Range("E5").Value = "S"
$ws.Range("F5").Value = "O"
$ws.Range("G5").Formula = "=_xlfn.CONCAT(""SFX_"",A5,""_"",B5,""_"",C5,""_"",D5)"
